# Ensure new comments are attributed to "Marc Smith" (matches the other
# comment already in this document).
$word.UserName = "Marc Smith"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "ardille" -> "ardile"
#    The word "ardille" appears five times in this document; only the
#    fourth occurrence (".../estends les sur l'ardille jusques a la
#    joincture...") is corrected by this edit, so walk the matches with
#    Find and patch the right one directly via a Range.
# ---------------------------------------------------------------------
$search = $d.Content
$hits = @()
while ($search.Find.Execute("ardille", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $hits += , @($search.Start, $search.End)
    $search.Collapse(0)
}

if ($hits.Count -ge 4) {
    $target = $hits[3]
    $fix = $d.Range($target[0], $target[1])
    $fix.Text = "ardile"
}

# ---------------------------------------------------------------------
# 2) Split "remuroient." after its first letter and anchor a comment
#    ("paper strip over the whole line") to that initial "r".
# ---------------------------------------------------------------------
$word_range = $d.Content
$word_range.Find.Execute("remuroient.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$commentAnchor = $d.Range($word_range.Start, $word_range.Start + 1)
$d.Comments.Add($commentAnchor, "paper strip over the whole line")
